$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "56.744.31"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.340.35"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  +0.06%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "514.74"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "133.89"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.534"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("E11").Value = "  +1.64%  "

$ws.Range("E12").Value = "  +0.18%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "23.87"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.759.84"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.72%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "56.703.34"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("E16").Value = "  +0.02%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.361.49"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.01%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "10.43"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "326.28"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.34%  "

$ws.Range("E20").Value = "  -1.05%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.69"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.18%  "

$ws.Range("E22").Value = "  +0.09%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "61.21"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.95%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "8.69"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +12.72%  "

$ws.Range("E25").Value = "  +4.01%  "

$ws.Range("E26").Value = "  +0.08%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.30"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +7.14%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "168.65"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("E30").Value = "  +0.68%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "6.15"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "18.42"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("E34").Value = "  -0.29%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.28"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +3.22%  "

$ws.Range("E36").Value = "  +0.88%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.890"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -5.58%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.57"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +3.12%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "38.62"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "150.83"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +9.48%  "

$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("E42").Value = "  +1.71%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "281.84"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.13"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +2.32%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0926"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0501"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.558"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "18.37"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +7.78%  "

$ws.Range("E49").Value = "  +0.27%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "17.13"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.48%  "

$ws.Range("E51").Value = "  +1.27%  "
